# cincy.reach.xlsx revision -- 23-Jun-17
# Adds the "bacterial growth efficiency" worksheet (whole-stream respiration /
# carbon-use efficiency calculations pulled from the "combined" sheet) and
# makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$combined = $wb.Worksheets.Item("combined")

# --- add the new worksheet as the last tab, after "combined" -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "bacterial growth efficiency"

# --- header row (written in the original authoring order so the shared-
#     string table comes out byte-identical to the source workbook) --------
$ws.Range("A1").Value = "stream"
$ws.Range("B1").Value = "season"
$ws.Range("C1").Value = "reach"
$ws.Range("D1").Value = "DOC (mg/L)"
$ws.Range("E1").Value = "Q (L/s)"
$ws.Range("G1").Value = "DOC (mg/s)"
$ws.Range("H1").Value = "Whole-stream Respiration (g O2 m-2 d-1)"
$ws.Range("I1").Value = "Whole-stream Respiration (g C m-2 d-1)"
$ws.Range("F1").Value = "w (m)"
$ws.Range("J1").Value = "Whole-stream Respiration (mg C m-2 s-1)"
$ws.Range("K1").Value = "sw.m"
$ws.Range("L1").Value = "k"

# --- data rows 2-18: pull from combined!, then derive the BGE metrics -----
# (Uses R1C1 relative formulas applied to the whole 2:18 block in one shot,
#  mirroring how the author filled these down from row 2 in Excel.)
$ws.Range("A2:A18").FormulaR1C1 = "=combined!RC[1]"
$ws.Range("B2:B18").FormulaR1C1 = "=combined!RC[1]"
$ws.Range("C2:C18").FormulaR1C1 = "=combined!RC[1]"
$ws.Range("D2:D18").FormulaR1C1 = "=combined!RC[46]"
$ws.Range("E2:E18").FormulaR1C1 = "=combined!RC[25]"
$ws.Range("G2:G18").FormulaR1C1 = "=RC[-2]*RC[-3]"
$ws.Range("H2:H18").FormulaR1C1 = "=combined!RC[27]"
$ws.Range("I2:I18").FormulaR1C1 = "=(RC[-1]/(15.999*2))*-12.011"
$ws.Range("F2:F18").FormulaR1C1 = "=combined!RC[27]"
$ws.Range("J2:J18").FormulaR1C1 = "=(RC[-1]*1000)/(24*60*60)"
$ws.Range("K2:K18").FormulaR1C1 = "=RC[-4]/(RC[-1]*RC[-5])"
$ws.Range("L2:L18").FormulaR1C1 = "=1/RC[-1]*-1"

# --- view: freeze header row off-screen at column Z, select AI2 -----------
$ws.Range("L2").Select()
$excel.ActiveWindow.FreezePanes = $false

# --- column widths (approximate Excel's auto-fit-to-content) --------------
$ws.Columns.Item(4).ColumnWidth = 11.140625
$ws.Columns.Item(7).ColumnWidth = 11.140625
$ws.Columns.Item(8).ColumnWidth = 38.140625
$ws.Columns.Item(9).ColumnWidth = 38.140625
$ws.Columns.Item(10).ColumnWidth = 36.7109375
$ws.Columns.Item(12).ColumnWidth = 12.5703125

# --- "combined" sheet view: scroll frozen pane to column Z, select AI2 ----
$combined.Activate()
$combined.Range("AI2").Select()
$excel.ActiveWindow.ScrollColumn = 26

# --- finally, re-activate the new sheet so it becomes the active tab ------
$ws.Activate()
$ws.Range("L2").Select()
